$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 20; this shifts rows 20-41 down to 21-42 and
# carries the formatting (styles) of the row above into the new row,
# matching the new blank-ish row seen in the diff (A20/E20 keep their
# style, B20 has no "ok" marker).
[void]$ws.Rows("20:20").Insert()

# Populate the two new "AdsList" (column G) entries for the new ads.
$ws.Cells.Item(20, 7).Value = "https://davao-city-das-ph.global-free-classified-ads.com/listings/camiguin-tour-package-philippines-it8993666.html"
$ws.Cells.Item(21, 7).Value = "https://davao-city-das-ph.global-free-classified-ads.com/listings/samal-island-hopping-it8993758.html"

# Match the author's final selection state.
[void]$ws.Range("H19").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 6
